$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function SetTextValue($Cell, $Text) {
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = $origStyle
}

SetTextValue $ws.Range('D2') '26.126.06'
$ws.Range('E2').Value = '  -0.20%  '
SetTextValue $ws.Range('D3') '1.654.90'
$ws.Range('E3').Value = '  -0.25%  '
$ws.Range('E4').Value = '  -0.33%  '
SetTextValue $ws.Range('D5') '217.86'
$ws.Range('E6').Value = '  +1.20%  '
$ws.Range('E7').Value = '  -0.26%  '
SetTextValue $ws.Range('D8') '0.2608'
$ws.Range('E8').Value = '  -0.77%  '
SetTextValue $ws.Range('D9') '0.06351'
$ws.Range('E9').Value = '  +1.79%  '
$ws.Range('E10').Value = '  -1.13%  '
SetTextValue $ws.Range('D11') '0.07783'
$ws.Range('E11').Value = '  +0.97%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
SetTextValue $ws.Range('D12') '4.501'
$ws.Range('E12').Value = '  +2.41%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
SetTextValue $ws.Range('D13') '1.621.32'
$ws.Range('E13').Value = '  -2.17%  '
SetTextValue $ws.Range('D14') '0.5484'
$ws.Range('E14').Value = '  +1.23%  '
SetTextValue $ws.Range('D15') '0.0₅8231'
$ws.Range('E15').Value = '  +1.94%  '
SetTextValue $ws.Range('D17') '26.126.91'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('E18').Value = '  -0.37%  '
SetTextValue $ws.Range('D19') '4.579'
$ws.Range('E19').Value = '  -0.58%  '
SetTextValue $ws.Range('D20') '191.60'
$ws.Range('E20').Value = '  +0.38%  '
$ws.Range('E21').Value = '  +0.63%  '
SetTextValue $ws.Range('D22') '6.043'
$ws.Range('E22').Value = '  +0.00%  '
SetTextValue $ws.Range('D23') '1.003'
$ws.Range('E23').Value = '  -0.39%  '
SetTextValue $ws.Range('D24') '142.09'
$ws.Range('E24').Value = '  +1.60%  '
SetTextValue $ws.Range('D25') '0.1250'
$ws.Range('E25').Value = '  +2.53%  '
$ws.Range('E26').Value = '  +2.01%  '
$ws.Range('E27').Value = '  +1.11%  '
SetTextValue $ws.Range('D28') '1.432'
$ws.Range('E28').Value = '  +1.61%  '
SetTextValue $ws.Range('D29') '0.05911'
$ws.Range('E29').Value = '  -1.00%  '
SetTextValue $ws.Range('D30') '1.282'
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('E31').Value = '  -0.98%  '
SetTextValue $ws.Range('D32') '3.255'
$ws.Range('E32').Value = '  +0.57%  '
SetTextValue $ws.Range('D33') '1.586'
$ws.Range('E33').Value = '  -1.61%  '
SetTextValue $ws.Range('D34') '0.9535'
$ws.Range('E34').Value = '  -0.67%  '
SetTextValue $ws.Range('D35') '2.782'
$ws.Range('E35').Value = '  +0.19%  '
SetTextValue $ws.Range('D36') '2.410'
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('E37').Value = '  +1.18%  '
SetTextValue $ws.Range('D38') '0.01620'
$ws.Range('E38').Value = '  +2.05%  '
SetTextValue $ws.Range('D39') '5.786'
$ws.Range('E39').Value = '  -3.67%  '
SetTextValue $ws.Range('D40') '0.8481'
$ws.Range('E40').Value = '  -0.90%  '
SetTextValue $ws.Range('D41') '1.002'
$ws.Range('E41').Value = '  -0.22%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
SetTextValue $ws.Range('D42') '103.08'
$ws.Range('E42').Value = '  +3.15%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
SetTextValue $ws.Range('D43') '1.029.29'
$ws.Range('E43').Value = '  +1.09%  '
SetTextValue $ws.Range('D44') '1.801.42'
$ws.Range('E44').Value = '  +0.08%  '
SetTextValue $ws.Range('D45') '57.28'
$ws.Range('E45').Value = '  +1.30%  '
SetTextValue $ws.Range('D46') '1.003'
$ws.Range('E46').Value = '  -0.87%  '
SetTextValue $ws.Range('D47') '0.4299'
$ws.Range('E48').Value = '  +1.75%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
SetTextValue $ws.Range('D49') '7.865'
$ws.Range('E49').Value = '  -0.65%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
SetTextValue $ws.Range('D50') '0.05163'
$ws.Range('E50').Value = '  -0.25%  '
SetTextValue $ws.Range('D51') '0.09721'
$ws.Range('E51').Value = '  +0.92%  '

Write-Output "Done applying changes"